# Adapt column header formatting to respective input file names (FV2310 / FV2404):
# the "_old"/"_new" header suffixes become "_FV2310"/"_FV2404", the data range is
# converted into a native Excel Table, and the header row is frozen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells: "<name>_old" -> "<name>_FV2310",
#        "<name>_new" -> "<name>_FV2404" -------------------------------------
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    # columns A..J (1..10) carry the "old" / FV2310 headers
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2310"
    # columns L..U (12..21) carry the "new" / FV2404 headers (column K stays "diff")
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2404"
}

# --- 2. Turn the used range into a native Excel Table -------------------------
$dataRange = $ws.UsedRange
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row --------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
